$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("table")

# Rename the header row from plural/grouped names to the new singular
# column names used by the table/range classes. Assign in the order
# name, odd, even, number so the shared-string table ends up ordered
# the same way the real edit produced it.
$ws.Range("E2").Value = "name"
$ws.Range("D2").Value = "odd"
$ws.Range("C2").Value = "even"
$ws.Range("B2").Value = "number"

# The original 10-row data block (rows 3-12) gets duplicated twice more,
# producing a 30-row table (rows 3-32).
$dataBlock = $ws.Range("B3:G12")
$dataBlock.Copy($ws.Range("B13"))
$dataBlock.Copy($ws.Range("B23"))

# The "tammy" record (6th record of each block) no longer carries a name
# value - clear the name cell in every copy of that record.
$ws.Range("E8").Value = ""
$ws.Range("E18").Value = ""
$ws.Range("E28").Value = ""

# Update the active selection to match where the editor ended up.
[void]$ws.Range("H15").Select()
